# Applies the Thu May 30 05:54:45 UTC 2024 "Updated cryptos list" GitHub Actions
# refresh: per-row Price (column D) / Volume(1h) (column E) updates, plus the two
# row swaps (Mantle<->Kaspa at rows 39-40, FLOKI<->Arweave at rows 50-51) where the
# Coin name (B) and Link (C) also change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "68.119.43"
$ws.Range("E2").Value = "  -0.81%  "
# Row 3 (Ethereum)
$ws.Range("D3").Value = "3.781.52"
$ws.Range("E3").Value = "  -2.35%  "
# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.00%  "
# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "
# Row 6 (Solana)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.61"
$ws.Range("D6").Style = "Normal"
# Row 7 (LidoStakedEther)
$ws.Range("D7").Value = "3.779.70"
$ws.Range("E7").Value = "  -2.40%  "
# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.03%  "
# Row 9 (XRP)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.67%  "
# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  -3.08%  "
# Row 11 (Toncoin)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.40%  "
# Row 12 (Cardano)
$ws.Range("E12").Value = "  -2.35%  "
# Row 13 (ShibaInu)
$ws.Range("E13").Value = "  -3.85%  "
# Row 14 (Avalanche)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.21%  "
# Row 15 (WrappedliquidstakedEther2.0)
$ws.Range("D15").Value = "4.416.40"
$ws.Range("E15").Value = "  -2.31%  "
# Row 16 (WrappedEther)
$ws.Range("D16").Value = "3.781.74"
$ws.Range("E16").Value = "  -2.53%  "
# Row 17 (Chainlink)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.65"
$ws.Range("D17").Style = "Normal"
# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "68.068.01"
$ws.Range("E18").Value = "  -1.04%  "
# Row 19 (Polkadot)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.75%  "
# Row 20 (TRON)
$ws.Range("E20").Value = "  -0.34%  "
# Row 21 (Uniswap)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.20%  "
# Row 22 (BitcoinCash)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.12%  "
# Row 23 (Polygon)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.717"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.30%  "
# Row 25 (Litecoin)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.24%  "
# Row 26 (Fetch.AI)
$ws.Range("E26").Value = "  -1.76%  "
# Row 27 (InternetComputer(DFINITY))
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.92%  "
# Row 28 (RenderToken)
$ws.Range("E28").Value = "  -0.56%  "
# Row 29 (Dai)
$ws.Range("E29").Value = "  -0.07%  "
# Row 30 (PancakeSwap)
$ws.Range("E30").Value = "  -1.18%  "
# Row 31 (WrappedeETH)
$ws.Range("D31").Value = "3.930.01"
$ws.Range("E31").Value = "  -2.34%  "
# Row 32 (NEARProtocol)
$ws.Range("E32").Value = "  -2.68%  "
# Row 33 (EthereumClassic)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.94%  "
# Row 34 (ImmutableX)
$ws.Range("E34").Value = "  -4.36%  "
# Row 35 ()
$ws.Range("E35").Value = "  -2.50%  "
# Row 36 (RenzoRestakedETH)
$ws.Range("D36").Value = "3.734.80"
$ws.Range("E36").Value = "  -2.71%  "
# Row 37 (Hedera)
$ws.Range("E37").Value = "  -1.78%  "
# Row 38 (dogwifhat)
$ws.Range("E38").Value = "  -6.48%  "
# Row 39 (Mantle -> Kaspa)
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.139"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.33%  "
# Row 40 (Kaspa -> Mantle)
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.56%  "
# Row 41 (Filecoin)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.81%  "
# Row 42 (FirstDigitalUSD)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
# Row 43 (TheGraph)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.313"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "
# Row 44 (USDe)
$ws.Range("E44").Value = "  +0.03%  "
# Row 46 (Stacks)
$ws.Range("E46").Value = "  -3.11%  "
# Row 47 (Bittensor)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "405.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.46%  "
# Row 48 (OKB)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.37%  "
# Row 49 (Monero)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "144.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.00%  "
# Row 50 (FLOKI -> Arweave)
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.30%  "
# Row 51 (Arweave -> FLOKI)
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000269"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.95%  "
